# Added Mounted Tire Processing Pipeline
# Updates row 6 ("signal segment 5") of the per-signal distribution (Step1_Data),
# its running cumulative total (Step2_Sj), and the threshold crossing values
# that feed the Step3_DataPts_* sheets.

$wb = $excel.ActiveWorkbook

# --- Step1_Data: row 6 normalized signal distribution --------------------
$ws1 = $wb.Worksheets.Item("Step1_Data")

$ws1.Range("D6").Value  = 0.007966266569831246
$ws1.Range("E6").Value  = 0.2493609080142864
$ws1.Range("F6").Value  = 0.08575276080584916
$ws1.Range("G6").Value  = 0.4197283025587936
$ws1.Range("H6").Value  = 0.01890901414044667
$ws1.Range("I6").Value  = 0.02809262200916679
$ws1.Range("J6").Value  = 0.02144000310803604
$ws1.Range("K6").Value  = 0.003406783932737229
$ws1.Range("O6").Value  = 0.003105022829482404
$ws1.Range("P6").Value  = 0.03747428763204592
$ws1.Range("Q6").Value  = 0.01123833775277237
$ws1.Range("S6").Value  = 0.00780747458089887
$ws1.Range("T6").Value  = 0.002565288839501491
$ws1.Range("U6").Value  = 0.01786009918713444
$ws1.Range("V6").Value  = 0.001403363510802767
$ws1.Range("Z6").Value  = 0.001863279539150565
$ws1.Range("AA6").Value = 0.006913617503949455
$ws1.Range("AC6").Value = 0.01860606903981144
$ws1.Range("AE6").Value = 0.0146443459926874
$ws1.Range("AF6").Value = 0.02060765981668817
$ws1.Range("AH6").Value = 0.01614195389236371
$ws1.Range("AI6").Value = 0.005112538743563869
$ws1.Range("AJ6").Value = 0

# --- Step2_Sj: row 6 cumulative sum of Step1_Data row 6 -------------------
$ws2 = $wb.Worksheets.Item("Step2_Sj")

$ws2.Range("D6").Value  = 0.007966266569831246
$ws2.Range("E6").Value  = 0.2573271745841177
$ws2.Range("F6").Value  = 0.3430799353899668
$ws2.Range("G6").Value  = 0.7628082379487604
$ws2.Range("H6").Value  = 0.7817172520892071
$ws2.Range("I6").Value  = 0.8098098740983739
$ws2.Range("J6").Value  = 0.83124987720641
$ws2.Range("K6").Value  = 0.8346566611391472
$ws2.Range("L6").Value  = 0.8346566611391472
$ws2.Range("M6").Value  = 0.8346566611391472
$ws2.Range("N6").Value  = 0.8346566611391472
$ws2.Range("O6").Value  = 0.8377616839686296
$ws2.Range("P6").Value  = 0.8752359716006756
$ws2.Range("Q6").Value  = 0.886474309353448
$ws2.Range("R6").Value  = 0.886474309353448
$ws2.Range("S6").Value  = 0.8942817839343469
$ws2.Range("T6").Value  = 0.8968470727738483
$ws2.Range("U6").Value  = 0.9147071719609828
$ws2.Range("V6").Value  = 0.9161105354717856
$ws2.Range("W6").Value  = 0.9161105354717856
$ws2.Range("X6").Value  = 0.9161105354717856
$ws2.Range("Y6").Value  = 0.9161105354717856
$ws2.Range("Z6").Value  = 0.9179738150109361
$ws2.Range("AA6").Value = 0.9248874325148856
$ws2.Range("AB6").Value = 0.9248874325148856
$ws2.Range("AC6").Value = 0.943493501554697
$ws2.Range("AD6").Value = 0.943493501554697
$ws2.Range("AE6").Value = 0.9581378475473844
$ws2.Range("AF6").Value = 0.9787455073640726
$ws2.Range("AG6").Value = 0.9787455073640726
$ws2.Range("AH6").Value = 0.9948874612564363
$ws2.Range("AI6").Value = 1

# --- Step3_DataPts_*: F6 holds the Step2_Sj value crossing that sheet's
#     intensity threshold (first Sj value > threshold, picked from row 6) --
$ws3 = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws3.Range("F6").Value = 0.7628082379487604

$ws4 = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws4.Range("F6").Value = 0.7628082379487604

$ws5 = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws5.Range("F6").Value = 0.8098098740983739

$ws6 = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws6.Range("F6").Value = 0.9147071719609828
